# "commit on 16th Final"
#
# Keyword_Scenario!C3 and Keyword_Scenario!C4:
#     "Application_Submit"  ->  "Application_Submit1"
# Test_Scenarios!D4:
#     "Yes"  ->  "No"
#
# The active sheet/selection state also moves: Keyword_Scenario becomes the
# selected tab with C4 as the active cell, while Test_Scenarios (previously
# the selected tab) keeps its D4 selection but is no longer active.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Keyword_Scenario")
$ws2 = $wb.Worksheets.Item("Test_Scenarios")

# -- data edits --
$ws1.Range("C3").Value = "Application_Submit1"
$ws1.Range("C4").Value = "Application_Submit1"
$ws2.Range("D4").Value = "No"

# -- view / selection state --
# Record Test_Scenarios' own selection (D4) without leaving it active.
$ws2.Range("D4").Select() | Out-Null

# Keyword_Scenario ends up the active sheet with C4 selected.
$ws1.Activate() | Out-Null
$ws1.Range("C4").Select() | Out-Null
